$d = $word.ActiveDocument
[void]$d.Content.Delete()
$r = $d.Content
$r.Collapse(1)
[void]$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:sz w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="32"/></w:rPr><w:t>Einleitung</w:t></w:r></w:p><w:p><w:r><w:t>Fragen/Teilgebiete/Gliederungspunkte/Absätze:</w:t></w:r></w:p><w:p><w:r><w:t>Motivation?</w:t></w:r></w:p><w:p><w:r><w:t>Aufmerksamkeit kann vereinfacht als begrenzte Ressource angesehen werden.</w:t></w:r></w:p><w:p><w:r><w:t>Wovon handelt die Arbeit?</w:t></w:r></w:p><w:p><w:r><w:t>Was ist ihr Ziel?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Welche </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Erkenntisse</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> sind zu finden?</w:t></w:r></w:p><w:p><w:r><w:t>Wie kann ich zu dem Thema hinführen?</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:sz w:val="32"/></w:rPr><w:t>Sinneswiederherstellung</w:t></w:r></w:p><w:p><w:r><w:t>Fragen/Teilgebiete/Gliederungspunkte/Absätze:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">eine Möglichkeit zur Kommunikation mit der Außenwelt. Im </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>folgenden</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> soll beschrieben werden, wie diese Eigenschaft genutzt werden kann um, über haptische Schnittstellen, ausgewählte Sinne wiederherstellen zu können.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:sz w:val="32"/></w:rPr><w:t>sehen</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Die Augen stellen eine mächtige Verbindung zur </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ausenwelt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> da. Deshalb ist eine eins </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>zu eins Übersetzung</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> über die Haut nur schwer vorstellbar. Deshalb geht es bei dieser Fragestellung darum die Komplexität der Informationen zu reduzieren. Beispielsweise könnte geschriebene Schrift von einer Kamera erfasst und in eine </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Brail</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">-artige Schrift übersetzt werden, die unter dem Finger des Anwenders manifestiert wird. Diese direkte Übersetzung bietet eine gute Möglichkeit das Prinzip der Komplexitätsreduktion zu erkennen. Das Problem dabei ist auch, dass die Haut nicht beliebig schnell Unterschiede wahrnehmen kann. Außerdem ist die Interpretation der Signale durch den Menschen ein weiterer Engpass. So wäre es </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ansonste</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> beispielsweise vorstellbar das Übersetzungsproblem durch 26 haptische Aktoren zu lösen. Dabei würde jeder Aktor zu einem Buchstaben im Alphabet zugeordnet werden. Die einzelnen Aktoren seien entlang des Unterarmes angeordnet. Die Differenzierung der Aktoren ist jetzt jedoch zu anspruchsvoll, wenn sich die Aktivierungsmuster der Aktoren nur durch ihre Position auf dem Körper unterscheiden.</w:t></w:r></w:p><w:p><w:r><w:t>\</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>subsection</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Leistungssteigerung}</w:t></w:r></w:p><w:p><w:r><w:t>Fragen/Teilgebiete/Gliederungspunkte/Absätze:</w:t></w:r></w:p><w:p><w:r><w:t>Leistung ist nach der Physik Arbeit pro Zeit. Um eine Leistungssteigerung zu erreichen muss also entweder die geleistete Arbeit bei gleicher Zeit erhöht werden oder dementsprechend die Zeit kürzer werden, die für eine Aufgabe gefragt ist.</w:t></w:r></w:p><w:p><w:r><w:t>Hier können haptische Aktuatoren unterstützend eingreifen.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>\</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>subsection</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Erweiterung des Wahrnehmungsspektrums}</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Fragen/Teilgebiete/Gliederungspunkte/Absätze:</w:t></w:r></w:p><w:p><w:r><w:t>Das Wahrnehmungsspektrum des Menschen ist durch die ihm zur Verfügung stehenden Sinne begrenzt. Auch spielt die Verarbeitungsgeschwindigkeit dieser Informationen für die Gesamtwahrnehmung eine Rolle.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Die evolutionäre Aufgabe der </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>haptik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ist auf kurzer Distanz (Berührung) Informationen über die Umwelt zu liefern. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dementsprechen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ist es von der Natur nicht vorgesehen größere Distanzen haptisch zu erfassen.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Jedoch sind Situationen denkbar, in den eine Verlagerung der Umgebungsanalyse von den Augen, die die Hauptverantwortlichen hierfür sind, auf andere Sinnesorgane </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vorzunehemen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. Dadurch werden andere Sinneskapazitäten freigeräumt.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')
